$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.964.55'
$ws.Range("E2").Value = '  -1.50%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.244.91'
$ws.Range("E3").Value = '  -1.70%  '

# Row 4
$ws.Range("E4").Value = '  +0.30%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.24'
$ws.Range("E5").Value = '  -0.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '298.74'
$ws.Range("E6").Value = '  +12.06%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  -2.00%  '

# Row 8
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +0.74%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.09'
$ws.Range("E10").Value = '  -3.78%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  -0.68%  '

# Row 12
$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.39'
$ws.Range("E12").Value = '  +3.17%  '

# Row 13
$ws.Range("E13").Value = '  -0.71%  '

# Row 14
$ws.Range("E14").Value = '  -2.72%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.34'
$ws.Range("E15").Value = '  -0.79%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.894'
$ws.Range("E16").Value = '  +1.88%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.581.64'
$ws.Range("E17").Value = '  -1.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.252.92'
$ws.Range("E18").Value = '  -1.32%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.845.04'
$ws.Range("E19").Value = '  -1.80%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.66'
$ws.Range("E20").Value = '  +11.46%  '

# Row 21
$ws.Range("E21").Value = '  -1.66%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.67'
$ws.Range("E22").Value = '  +26.82%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.60'
$ws.Range("E23").Value = '  +1.51%  '

# Row 24
$ws.Range("E24").Value = '  -3.85%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '232.12'
$ws.Range("E25").Value = '  -1.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.49'
$ws.Range("E26").Value = '  -0.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  +4.13%  '

# Row 28
$ws.Range("E28").Value = '  -1.60%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.18'
$ws.Range("E29").Value = '  -4.87%  '

# Row 30
$ws.Range("E30").Value = '  -1.07%  '

# Row 31
$ws.Range("E31").Value = '  -4.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '175.42'
$ws.Range("E32").Value = '  +0.92%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.30'
$ws.Range("E33").Value = '  -1.84%  '

# Row 34
$ws.Range("E34").Value = '  -0.75%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.70'
$ws.Range("E35").Value = '  -0.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.51'
$ws.Range("E36").Value = '  +14.35%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.129'
$ws.Range("E37").Value = '  -1.28%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.82'
$ws.Range("E38").Value = '  +2.92%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0372'
$ws.Range("E39").Value = '  -3.52%  '

# Row 40
$ws.Range("E40").Value = '  -1.15%  '

# Row 41
$ws.Range("E41").Value = '  +1.50%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.241'
$ws.Range("E42").Value = '  +2.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.98'
$ws.Range("E43").Value = '  -3.11%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.42'
$ws.Range("E44").Value = '  -5.23%  '

# Row 45
$ws.Range("E45").Value = '  +0.21%  '

# Row 46
$ws.Range("E46").Value = '  -2.41%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.59'
$ws.Range("E47").Value = '  -6.32%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.37'
$ws.Range("E48").Value = '  +7.67%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.92'
$ws.Range("E49").Value = '  +5.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.67'
$ws.Range("E50").Value = '  +0.73%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0989'
$ws.Range("E51").Value = '  -1.52%  '
